$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$cells = @("B65", "B75", "B79", "B85", "B86", "B96", "B99", "B100", "B106", "B107", "B108", "B109")

foreach ($addr in $cells) {
    $rng = $ws.Range($addr)
    $orig = $rng.Value2
    $rng.Value = $orig -replace '\|', '||'
}

$ws.Activate()
$ws.Range("B109").Select()
